$wb = $excel.ActiveWorkbook

# --- Sheet "Gaz" ---
$gaz = $wb.Worksheets.Item("Gaz")

# Update header row (row 1): columns B/C renamed, D now holds what used to be
# in the "Last Price"/"Last Volume"/"End of Day Index" columns (E:G), so the
# unused E:G columns are cleared away entirely (content + formatting).
$gaz.Range("B1").Value = "Last Price"
$gaz.Range("C1").Value = "Last Volume"
$gaz.Range("D1").Value = "End of Day Index"

# Row 2 (2025-06-16): bid/ask/last placeholders become real figures
$gaz.Range("B2").Value = 37.15
$gaz.Range("C2").Value = 13680
$gaz.Range("D2").Value = 36.934

# Row 3 (2025-06-17): values that used to live in E3:G3 move into B3:D3
$gaz.Range("B3").Value = 38.95
$gaz.Range("C3").Value = 24000
$gaz.Range("D3").Value = 38.201

# Remove the now-unused E:G columns (rows 1-3) entirely so the sheet
# dimension shrinks back down to A1:D3
$gaz.Range("E1:G3").Clear()

# --- Sheet "CO2" ---
$co2 = $wb.Worksheets.Item("CO2")
$co2.Range("B2").Value = 74.7
